{"js": "// The document contains a practice sheet (\"two-digit number divided by\n// one-digit number\") laid out as a table of \"NN\u00f7N=\" problems. This edit\n// swaps out the 25 problem strings (in document order) for a new set of\n// problems, leaving the date heading, table layout and formatting intact.\n//\n// Because several of the original problem strings repeat (e.g. \"36\u00f77=\"\n// appears three times but maps to three different replacements), a plain\n// global find/replace is not safe here \u2014 each occurrence must be updated\n// independently, in the order the problems appear in the document.\n\nconst replacements = [\n  [\"97\u00f78=\", \"41\u00f76=\"],\n  [\"10\u00f76=\", \"83\u00f76=\"],\n  [\"11\u00f74=\", \"91\u00f76=\"],\n  [\"98\u00f73=\", \"66\u00f77=\"],\n  [\"36\u00f77=\", \"38\u00f75=\"],\n  [\"25\u00f75=\", \"32\u00f72=\"],\n  [\"61\u00f77=\", \"18\u00f79=\"],\n  [\"21\u00f73=\", \"40\u00f76=\"],\n  [\"22\u00f75=\", \"45\u00f72=\"],\n  [\"25\u00f79=\", \"95\u00f79=\"],\n  [\"36\u00f77=\", \"40\u00f75=\"],\n  [\"19\u00f73=\", \"88\u00f76=\"],\n  [\"59\u00f75=\", \"79\u00f76=\"],\n  [\"89\u00f75=\", \"59\u00f75=\"],\n  [\"53\u00f79=\", \"82\u00f79=\"],\n  [\"36\u00f77=\", \"50\u00f75=\"],\n  [\"77\u00f73=\", \"65\u00f74=\"],\n  [\"97\u00f72=\", \"82\u00f72=\"],\n  [\"44\u00f78=\", \"56\u00f74=\"],\n  [\"15\u00f74=\", \"84\u00f77=\"],\n  [\"89\u00f73=\", \"30\u00f75=\"],\n  [\"35\u00f72=\", \"18\u00f72=\"],\n  [\"17\u00f79=\", \"84\u00f77=\"],\n  [\"50\u00f72=\", \"70\u00f74=\"],\n  [\"46\u00f73=\", \"57\u00f76=\"],\n];\n\n// Matches the \"two digit/one digit \u00f7 problem\" text, e.g. \"97\u00f78=\".\nconst problemPattern = /^\\s*\\d+\u00f7\\d+=\\s*$/;\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet replacementIndex = 0;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (replacementIndex >= replacements.length) {\n    break;\n  }\n  const paragraph = paragraphs.items[i];\n  const text = paragraph.text;\n  if (!problemPattern.test(text)) {\n    continue;\n  }\n  const [expected, updated] = replacements[replacementIndex];\n  if (text !== expected) {\n    throw new Error(\n      `Unexpected problem text at position ${replacementIndex}: ` +\n        `found \"${text}\", expected \"${expected}\"`\n    );\n  }\n  paragraph.insertText(updated, Word.InsertLocation.replace);\n  replacementIndex++;\n}\n\nawait context.sync();\n\nif (replacementIndex !== replacements.length) {\n  throw new Error(\n    `Only replaced ${replacementIndex} of ${replacements.length} problems`\n  );\n}\n", "ps1": "# The document is a practice sheet (\"two-digit number divided by\n# one-digit number\") laid out as a table of \"NN\u00f7N=\" problems. This edit\n# swaps out the 25 problem strings (in document order) for a new set of\n# problems, leaving the date heading, table layout and formatting intact.\n#\n# Because several of the original problem strings repeat (e.g. \"36\u00f77=\"\n# appears three times but maps to three different replacements), a plain\n# global Find & Replace is not safe here -- each occurrence must be\n# updated independently, in the order the problems appear in the\n# document.\n\n$replacements = @(\n    ,@(\"97\u00f78=\", \"41\u00f76=\")\n    ,@(\"10\u00f76=\", \"83\u00f76=\")\n    ,@(\"11\u00f74=\", \"91\u00f76=\")\n    ,@(\"98\u00f73=\", \"66\u00f77=\")\n    ,@(\"36\u00f77=\", \"38\u00f75=\")\n    ,@(\"25\u00f75=\", \"32\u00f72=\")\n    ,@(\"61\u00f77=\", \"18\u00f79=\")\n    ,@(\"21\u00f73=\", \"40\u00f76=\")\n    ,@(\"22\u00f75=\", \"45\u00f72=\")\n    ,@(\"25\u00f79=\", \"95\u00f79=\")\n    ,@(\"36\u00f77=\", \"40\u00f75=\")\n    ,@(\"19\u00f73=\", \"88\u00f76=\")\n    ,@(\"59\u00f75=\", \"79\u00f76=\")\n    ,@(\"89\u00f75=\", \"59\u00f75=\")\n    ,@(\"53\u00f79=\", \"82\u00f79=\")\n    ,@(\"36\u00f77=\", \"50\u00f75=\")\n    ,@(\"77\u00f73=\", \"65\u00f74=\")\n    ,@(\"97\u00f72=\", \"82\u00f72=\")\n    ,@(\"44\u00f78=\", \"56\u00f74=\")\n    ,@(\"15\u00f74=\", \"84\u00f77=\")\n    ,@(\"89\u00f73=\", \"30\u00f75=\")\n    ,@(\"35\u00f72=\", \"18\u00f72=\")\n    ,@(\"17\u00f79=\", \"84\u00f77=\")\n    ,@(\"50\u00f72=\", \"70\u00f74=\")\n    ,@(\"46\u00f73=\", \"57\u00f76=\")\n)\n\n$d = $word.ActiveDocument\n\n$replacementIndex = 0\nforeach ($p in $d.Paragraphs) {\n    if ($replacementIndex -ge $replacements.Count) {\n        break\n    }\n\n    $r = $p.Range\n    $rawText = $r.Text\n    # Strip the trailing paragraph mark / table-cell mark control chars.\n    $cleanText = $rawText -replace '[\\x00-\\x1F\\x7F]+$', ''\n\n    if ($cleanText -notmatch '^\\d+\u00f7\\d+=$') {\n        continue\n    }\n\n    $pair = $replacements[$replacementIndex]\n    $expected = $pair[0]\n    $updated = $pair[1]\n\n    if ($cleanText -ne $expected) {\n        throw \"Unexpected problem text at position $replacementIndex`: found [$cleanText], expected [$expected]\"\n    }\n\n    # Replace only the text portion, leaving the trailing marks (and run\n    # formatting) untouched.\n    $textRange = $r.Duplicate\n    $textRange.End = $textRange.Start + $cleanText.Length\n    $textRange.Text = $updated\n\n    $replacementIndex++\n}\n\nif ($replacementIndex -ne $replacements.Count) {\n    throw \"Only replaced $replacementIndex of $($replacements.Count) problems\"\n}\n"}
